# Aggiornato ObjectID a MagazzionoF
# The "objMapping" sheet (sheet2) encodes, per group, a header row
# ("viewXxx") followed by its "xxxChild1..5" rows. The "viewMagazzino"
# group is shrunk from 5 children to 3, and a brand-new "viewMagazzinoF"
# group (5 children) is appended right before the closing "]" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Fix typo in the opening label (row 1, col F) ---
$ws.Range("F1").Value = "objMapping = ["

# --- Column A: shift rows 19-38 up by two (viewMagazzino now only has
#     magazzinoChild1-3, so everything from viewPreparazione onward moves
#     up two rows), then append the six new viewMagazzinoF rows (37-42).
$newA = @{
    19 = "viewPreparazione";    20 = "preparazioneChild1"
    21 = "preparazioneChild2";  22 = "preparazioneChild3"
    23 = "preparazioneChild4";  24 = "preparazioneChild5"
    25 = "viewLavorazione";     26 = "lavorazioneChild1"
    27 = "lavorazioneChild2";   28 = "lavorazioneChild3"
    29 = "lavorazioneChild4";   30 = "lavorazioneChild5"
    31 = "viewFinitura";        32 = "finituraChild1"
    33 = "finituraChild2";      34 = "finituraChild3"
    35 = "finituraChild4";      36 = "finituraChild5"
    37 = "viewMagazzinoF";      38 = "magazzinoFChild1"
    39 = "magazzinoFChild2";    40 = "magazzinoFChild3"
    41 = "magazzinoFChild4";    42 = "magazzinoFChild5"
}
foreach ($r in 19..42) {
    $ws.Range("A$r").Value = $newA[$r]
}

# --- C11/C12/C13: the "chartN shows header of group N" rows point at
#     the (now two-rows-earlier) viewPreparazione/viewLavorazione/
#     viewFinitura header rows. C14 now mirrors the new viewMagazzinoF
#     group header (row 37) instead of the viewPiantina header (row 9). ---
$ws.Range("C11").Formula = "=A19"
$ws.Range("C12").Formula = "=A25"
$ws.Range("C13").Formula = "=A31"
$ws.Range("C14").Formula = "=A37"

# --- Extend columns B, C, E, F for the six new rows (37-42) so the table
#     keeps the same per-row pattern used by every other group. ---
foreach ($r in 37..42) {
    $ws.Range("B$r").Value = ","
    $ws.Range("C$r").Formula = "=A`$9"
    if ($r -gt 37) {
        $ws.Range("E$r").Formula = "=E" + ($r - 1) + "+1"
    }
    $ws.Range("F$r").Formula = '="    ObjectID."&C' + $r + ' &","'
}

# --- Move the closing "]" marker from F39 down to F43 (the new last row),
#     and populate the new F39:F42 cells with the per-row formula. ---
$ws.Range("F43").Value = "]"

# --- Selection, to match the saved workbook view. ---
$ws.Range("C15").Select()
